$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New TPM-derived values for columns E:T, rows 2-13
$data = @{
    2 = @(1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 2, 1, 2.4249555, 4.849911, 0.1414947904552067, 0.1135050984851495, 1.7578720665495, 7.031488266198, 0.08184355702079261, 0.05422390218038106)
    3 = @(1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 3, 1, 7.537927333333333, 22.613782, 0.4398338231347236, 0.5292426094069979, 5.464311365279333, 32.785868191676, 0.2544091161772898, 0.2528309288761098)
    4 = @(1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 1, 0.3333333333333333, 0.1035233333333333, 0.31057, 0.006040528313704939, 0.007268438211862629, 0.07504499604333334, 0.45026997626, 0.003493968377831753, 0.003472294089553591)
    5 = @(1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 2, 1, 6.260847500000001, 12.521695, 0.3653169326548487, 0.2930520218156591, 4.538544700377501, 18.15417880151, 0.2113069824847247, 0.1399974483681385)
    6 = @(1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 3, 1, 0.6812176666666666, 2.043653, 0.03974866796499352, 0.04782871995681392, 0.4938208175256666, 2.962924905154, 0.02299146394455677, 0.02284883998131972)
    7 = @(1, 0.5, 0.724909, 1.449818, 0.5784209917375155, 0.477722171991027, 2, 0.6666666666666666, 0.1296543333333333, 0.388963, 0.007565257476522569, 0.009103112123517158, 0.09398759312233335, 0.5639255587340001, 0.004375903732319838, 0.004348758495524466)
    8 = @(2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 2, 1, 2.4249555, 4.849911, 0.1414947904552067, 0.1135050984851495, 1.281215538603, 7.687293231617999, 0.0596512334344141, 0.05928119630476848)
    9 = @(2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 3, 1, 7.537927333333333, 22.613782, 0.4398338231347236, 0.5292426094069979, 3.982633754857333, 35.843703793716, 0.1854247069574338, 0.2764116805308881)
    10 = @(2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 1, 0.3333333333333333, 0.1035233333333333, 0.31057, 0.006040528313704939, 0.007268438211862629, 0.05469613907333334, 0.49226525166, 0.002546559935873186, 0.003796144122309038)
    11 = @(2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 2, 1, 6.260847500000001, 12.521695, 0.3653169326548487, 0.2930520218156591, 3.307893733235, 19.84736239941, 0.1540099501701239, 0.1530545734475206)
    12 = @(2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 3, 1, 0.6812176666666666, 2.043653, 0.03974866796499352, 0.04782871995681392, 0.3599186293126666, 3.239267663814, 0.01675720402043676, 0.0249798799754942)
    13 = @(2, 0.6666666666666666, 0.528346, 1.585038, 0.4215790082624845, 0.522277828008973, 2, 0.6666666666666666, 0.1296543333333333, 0.388963, 0.007565257476522569, 0.009103112123517158, 0.06850234839933333, 0.616521135594, 0.00318935374420273, 0.004754353627992692)
}

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
